$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System'
$ws.Range("G3").Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg'
$ws.Range("G4").Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G5").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G6").Value = 'majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range("G7").Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G9").Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G11").Value = 'aya.saeed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G12").Value = 'Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg'
$ws.Range("G13").Value = 'yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G24").Value = 'youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range("G27").Value = 'nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg'
$ws.Range("G28").Value = 'maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G30").Value = 'shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
